{"js": "// Update the worksheet date and every \"A\u00d7B=C\" equation to the new set of\n// values (see commit \"Update master to output generated at 503736d\").\n// Each entry is [oldText, newText]; old values are unique in the document,\n// so a case-sensitive exact search safely locates the single run to replace.\nconst replacements = [\n  [\"2025-04-14 Monday\", \"2025-04-15 Tuesday\"],\n  [\"870\u00d76=5220\", \"262\u00d75=1310\"],\n  [\"237\u00d76=1422\", \"771\u00d75=3855\"],\n  [\"356\u00d78=2848\", \"614\u00d76=3684\"],\n  [\"995\u00d72=1990\", \"138\u00d77=966\"],\n  [\"534\u00d75=2670\", \"452\u00d77=3164\"],\n  [\"329\u00d77=2303\", \"806\u00d72=1612\"],\n  [\"196\u00d75=980\", \"188\u00d78=1504\"],\n  [\"949\u00d79=8541\", \"394\u00d77=2758\"],\n  [\"921\u00d74=3684\", \"765\u00d73=2295\"],\n  [\"459\u00d72=918\", \"832\u00d77=5824\"],\n  [\"242\u00d79=2178\", \"502\u00d75=2510\"],\n  [\"749\u00d72=1498\", \"553\u00d72=1106\"],\n  [\"692\u00d77=4844\", \"782\u00d75=3910\"],\n  [\"707\u00d78=5656\", \"718\u00d73=2154\"],\n  [\"315\u00d78=2520\", \"923\u00d78=7384\"],\n  [\"515\u00d76=3090\", \"892\u00d77=6244\"],\n  [\"398\u00d73=1194\", \"548\u00d77=3836\"],\n  [\"325\u00d77=2275\", \"975\u00d78=7800\"],\n  [\"465\u00d72=930\", \"929\u00d78=7432\"],\n  [\"353\u00d78=2824\", \"641\u00d75=3205\"],\n  [\"647\u00d72=1294\", \"325\u00d76=1950\"],\n  [\"963\u00d72=1926\", \"336\u00d74=1344\"],\n  [\"422\u00d74=1688\", \"776\u00d77=5432\"],\n  [\"645\u00d79=5805\", \"586\u00d76=3516\"],\n  [\"362\u00d76=2172\", \"489\u00d75=2445\"],\n];\n\nconst body = context.document.body;\nlet replacedCount = 0;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n    replacedCount++;\n  }\n}\n\nawait context.sync();\nreturn \"replaced=\" + replacedCount;\n", "ps1": "# Update the worksheet date and every \"A x B = C\" equation to the new set\n# of values (see commit \"Update master to output generated at 503736d\").\n# $oldTexts[i] -> $newTexts[i]; old values are unique in the document, so\n# an exact, case-sensitive Find/Replace safely targets a single run each.\n$d = $word.ActiveDocument\n\n$oldTexts = @(\n    \"2025-04-14 Monday\",\n    \"870\u00d76=5220\",\n    \"237\u00d76=1422\",\n    \"356\u00d78=2848\",\n    \"995\u00d72=1990\",\n    \"534\u00d75=2670\",\n    \"329\u00d77=2303\",\n    \"196\u00d75=980\",\n    \"949\u00d79=8541\",\n    \"921\u00d74=3684\",\n    \"459\u00d72=918\",\n    \"242\u00d79=2178\",\n    \"749\u00d72=1498\",\n    \"692\u00d77=4844\",\n    \"707\u00d78=5656\",\n    \"315\u00d78=2520\",\n    \"515\u00d76=3090\",\n    \"398\u00d73=1194\",\n    \"325\u00d77=2275\",\n    \"465\u00d72=930\",\n    \"353\u00d78=2824\",\n    \"647\u00d72=1294\",\n    \"963\u00d72=1926\",\n    \"422\u00d74=1688\",\n    \"645\u00d79=5805\",\n    \"362\u00d76=2172\"\n)\n$newTexts = @(\n    \"2025-04-15 Tuesday\",\n    \"262\u00d75=1310\",\n    \"771\u00d75=3855\",\n    \"614\u00d76=3684\",\n    \"138\u00d77=966\",\n    \"452\u00d77=3164\",\n    \"806\u00d72=1612\",\n    \"188\u00d78=1504\",\n    \"394\u00d77=2758\",\n    \"765\u00d73=2295\",\n    \"832\u00d77=5824\",\n    \"502\u00d75=2510\",\n    \"553\u00d72=1106\",\n    \"782\u00d75=3910\",\n    \"718\u00d73=2154\",\n    \"923\u00d78=7384\",\n    \"892\u00d77=6244\",\n    \"548\u00d77=3836\",\n    \"975\u00d78=7800\",\n    \"929\u00d78=7432\",\n    \"641\u00d75=3205\",\n    \"325\u00d76=1950\",\n    \"336\u00d74=1344\",\n    \"776\u00d77=5432\",\n    \"586\u00d76=3516\",\n    \"489\u00d75=2445\"\n)\n\nfor ($i = 0; $i -lt $oldTexts.Count; $i++) {\n    $oldText = $oldTexts[$i]\n    $newText = $newTexts[$i]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # wdFindContinue=1 (Wrap), wdReplaceAll=2 (ReplaceWith) -- each $oldText is\n    # unique in the document, so ReplaceAll touches exactly one occurrence.\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
